$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price cells so values like "9.00" or "0.108"
# are preserved exactly as text instead of being auto-converted to numbers by Excel.

$ws.Range("D2").Value = '44.216.23'
$ws.Range("E2").Value = '  +2.96%  '

$ws.Range("D3").Value = '2.286.74'
$ws.Range("E3").Value = '  +2.85%  '

$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.38'
$ws.Range("E5").Value = '  +0.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.36'
$ws.Range("E6").Value = '  +4.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.588'
$ws.Range("E7").Value = '  +1.59%  '

$ws.Range("E8").Value = '  -0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.574'
$ws.Range("E9").Value = '  +1.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.88'
$ws.Range("E10").Value = '  +6.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0841'
$ws.Range("E11").Value = '  +2.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.89'
$ws.Range("E12").Value = '  +2.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("E13").Value = '  +2.16%  '

$ws.Range("D14").Value = '2.629.88'
$ws.Range("E14").Value = '  +2.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.880'
$ws.Range("E15").Value = '  +2.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.61'
$ws.Range("E16").Value = '  +4.24%  '

$ws.Range("D17").Value = '2.288.94'
$ws.Range("E17").Value = '  +3.24%  '

$ws.Range("D18").Value = '44.107.85'
$ws.Range("E18").Value = '  +2.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.49'
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").Value = '0.0₃0995'
$ws.Range("E20").Value = '  +3.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("E21").Value = '  +3.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.26'
$ws.Range("E22").Value = '  +1.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.23'
$ws.Range("E23").Value = '  +2.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '238.94'
$ws.Range("E24").Value = '  +1.33%  '

$ws.Range("E25").Value = '  +2.92%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.28'
$ws.Range("E27").Value = '  +1.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.89'
$ws.Range("E28").Value = '  +15.14%  '

$ws.Range("E29").Value = '  +1.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.56'
$ws.Range("E30").Value = '  +4.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.53'
$ws.Range("E31").Value = '  +3.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.55'
$ws.Range("E32").Value = '  +0.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0884'
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.72'
$ws.Range("E34").Value = '  -1.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.32'
$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.08'
$ws.Range("E36").Value = '  +3.30%  '

$ws.Range("E37").Value = '  -0.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.57'
$ws.Range("E38").Value = '  +2.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.108'
$ws.Range("E39").Value = '  +4.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.90'
$ws.Range("E40").Value = '  +6.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.71'
$ws.Range("E41").Value = '  +29.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0328'
$ws.Range("E42").Value = '  +1.58%  '

$ws.Range("E43").Value = '  -0.23%  '

$ws.Range("D44").Value = '1.783.18'
$ws.Range("E44").Value = '  -3.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.208'
$ws.Range("E45").Value = '  +1.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '85.63'
$ws.Range("E46").Value = '  -2.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.44'
$ws.Range("E47").Value = '  -0.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.00'
$ws.Range("E48").Value = '  +3.81%  '

$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.26'
$ws.Range("E49").Value = '  -4.04%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '59.81'
$ws.Range("E50").Value = '  -0.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '105.18'
$ws.Range("E51").Value = '  +4.61%  '
